{"js": "// Insert two new \"widget property\" list items after the\n// \"add space between 2 widgets.\" bullet and before the\n// \"Flutter Layout Cheat Sheet:\" bullet, then merge the\n// \"Flutter Layout \" / \"Cheat Sheet: \" runs into a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"add space between 2 widgets.\" bullet (ilvl=3) which\n// immediately precedes the \"Flutter Layout Cheat Sheet:\" bullet.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"add space between 2 widgets.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find anchor paragraph \"add space between 2 widgets.\"');\n}\n\n// New bullet: \"CircleAvatar\" at list level 2 (same level as \"SizedBox\").\nconst circleAvatarPara = anchor.insertParagraph(\"CircleAvatar\", \"After\");\ncircleAvatarPara.listItemOrNullObject.level = 2;\n\n// New bullet: \"backgroundImage: set image background.\" at list level 3\n// (same level as \"add space between 2 widgets.\").\nconst backgroundImagePara = circleAvatarPara.insertParagraph(\n  \"backgroundImage\",\n  \"After\"\n);\nbackgroundImagePara.listItemOrNullObject.level = 3;\nbackgroundImagePara\n  .getRange(\"End\")\n  .insertText(\": set image background.\", \"End\");\n\nawait context.sync();\n\n// Merge the \"Flutter Layout \" and \"Cheat Sheet: \" runs into one run\n// reading \"Flutter Layout Cheat Sheet: \" (they used to be split by the\n// _GoBack bookmark).\nconst flutterLayoutResults = body.search(\"Flutter Layout \", {\n  matchCase: true,\n});\nconst cheatSheetResults = body.search(\"Cheat Sheet: \", { matchCase: true });\nawait context.sync();\n\nif (flutterLayoutResults.items.length && cheatSheetResults.items.length) {\n  const merged = flutterLayoutResults.items[0].expandTo(\n    cheatSheetResults.items[0]\n  );\n  merged.insertText(\"Flutter Layout Cheat Sheet: \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Insert two new \"widget property\" list items after the\n# \"add space between 2 widgets.\" bullet and before the\n# \"Flutter Layout Cheat Sheet:\" bullet, then merge the\n# \"Flutter Layout \" / \"Cheat Sheet: \" runs into a single run.\n\n$doc = $word.ActiveDocument\n\n# Locate the \"add space between 2 widgets.\" bullet (list level 4, i.e.\n# w:ilvl=3) which immediately precedes the \"Flutter Layout Cheat Sheet:\"\n# bullet.\n$anchorIndex = -1\n$idx = 1\nforeach ($p in $doc.Paragraphs) {\n    if ($p.Range.Text -like \"*add space between 2 widgets.*\") {\n        $anchorIndex = $idx\n        break\n    }\n    $idx = $idx + 1\n}\nif ($anchorIndex -eq -1) {\n    throw 'Could not find anchor paragraph \"add space between 2 widgets.\"'\n}\n\n$anchor = $doc.Paragraphs.Item($anchorIndex)\n$anchorRange = $anchor.Range\n$anchorRange.Collapse(0)\n$anchorRange.InsertParagraphAfter()\n\n# New bullet: \"CircleAvatar\" at list level 3 (same level as \"SizedBox\",\n# i.e. w:ilvl=2).\n$circleAvatarPara = $doc.Paragraphs.Item($anchorIndex + 1)\n$circleAvatarPara.Range.Text = \"CircleAvatar\"\n$circleAvatarPara.Range.ListFormat.ListLevelNumber = 3\n\n$circleAvatarRange = $circleAvatarPara.Range\n$circleAvatarRange.Collapse(0)\n$circleAvatarRange.InsertParagraphAfter()\n\n# New bullet: \"backgroundImage: set image background.\" at list level 4\n# (same level as \"add space between 2 widgets.\", i.e. w:ilvl=3).\n$backgroundImagePara = $doc.Paragraphs.Item($anchorIndex + 2)\n$backgroundImagePara.Range.Text = \"backgroundImage: set image background.\"\n$backgroundImagePara.Range.ListFormat.ListLevelNumber = 4\n\n# Merge the \"Flutter Layout \" and \"Cheat Sheet: \" runs into one run\n# reading \"Flutter Layout Cheat Sheet: \" (they used to be split by the\n# _GoBack bookmark). Replace through a temporary placeholder first so the\n# engine performs a genuine text substitution (a same-text \"replace\" is a\n# no-op and would leave the original run split in place).\n$tmpMarker = \"FLCS_MERGE_TMP_0f3a\"\n\n$r1 = $doc.Range()\n$r1.Find.Execute(\"Flutter Layout Cheat Sheet: \", $false, $false, $false, $false, $false, $true, 1, $false, $tmpMarker, 2) | Out-Null\n\n$r2 = $doc.Range()\n$r2.Find.Execute($tmpMarker, $false, $false, $false, $false, $false, $true, 1, $false, \"Flutter Layout Cheat Sheet: \", 2) | Out-Null\n"}
